$d = $word.ActiveDocument

# 1) Replace "admin(owner)" with "admin(manager)" in the summary paragraph
$d.Content.Find.Execute("admin(owner)", $false, $false, $false, $false, $false, $true, 1, $false, "admin(manager)", 2)

# 2) Replace "Admin (Owner) Account Details" with "Admin (Manager) Account Details"
$d.Content.Find.Execute("Admin (Owner) Account Details", $false, $false, $false, $false, $false, $true, 1, $false, "Admin (Manager) Account Details", 2)
